$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 56
$ws.Range("F3").Value = 991
$ws.Range("F5").Value = 444
$ws.Range("F6").Value = 694
$ws.Range("F7").Value = 246
$ws.Range("F9").Value = 22
$ws.Range("F10").Value = 390
$ws.Range("F11").Value = 195
$ws.Range("F13").Value = 797
$ws.Range("F15").Value = 1961
$ws.Range("F16").Value = 459
$ws.Range("F17").Value = 6733
$ws.Range("F18").Value = 509
$ws.Range("F19").Value = 516
$ws.Range("F20").Value = 50
$ws.Range("F21").Value = 86
$ws.Range("F23").Value = 206

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 514
$ws.Range("F8").Value = 2
$ws.Range("F13").Value = 51

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5451
$ws.Range("F3").Value = 384
$ws.Range("F4").Value = 377

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 56
$ws.Range("F3").Value = 5451
$ws.Range("F4").Value = 384
$ws.Range("F5").Value = 377
$ws.Range("F10").Value = 514
$ws.Range("F11").Value = 991
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = 444
$ws.Range("F16").Value = 694
$ws.Range("F17").Value = 246
$ws.Range("F20").Value = 22
$ws.Range("F21").Value = 390
$ws.Range("F22").Value = 195
$ws.Range("F26").Value = 797
$ws.Range("F29").Value = 1961
$ws.Range("F30").Value = 459
$ws.Range("F31").Value = 6733
$ws.Range("F32").Value = 51
$ws.Range("F33").Value = 509
$ws.Range("F34").Value = 516
$ws.Range("F35").Value = 50
$ws.Range("F36").Value = 86
$ws.Range("F39").Value = 206
